# Update the "dSF" (column F) values for a subset of rows in Sheet1.
# These correspond to a data repull / recalculation of the mean (per the
# commit message "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -4
    "F3"  = -1
    "F4"  = -2
    "F6"  = -2
    "F7"  = 3
    "F9"  = 1
    "F11" = 5
    "F12" = 3
    "F16" = 2
    "F17" = 4
    "F18" = -3
    "F20" = 1
    "F21" = 8
    "F22" = -10
    "F23" = -8
    "F24" = -9
    "F25" = -9
    "F26" = -4
    "F28" = -7
    "F36" = -3
    "F40" = -9
    "F42" = -3
    "F43" = -3
    "F44" = -8
    "F46" = -1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
